$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a "numeric-looking" value as TEXT (leading apostrophe forces
# Excel to store it as a text/string cell instead of a number), matching the
# existing sheet convention where columns A/B/D hold textual values such as
# "103" and "123" rather than real numbers.
function Set-TextValue($row, $col, $text) {
    $ws.Cells.Item($row, $col).Value = "'" + $text
}

# --- New data rows (A, B, D populated as text; column C intentionally left
#     untouched/absent for rows 10, 16, 17, 18) ---

# Row 10
Set-TextValue 10 1 "123"
Set-TextValue 10 2 "123"
Set-TextValue 10 4 "123"

# Row 16
Set-TextValue 16 1 "123"
Set-TextValue 16 2 "123"
Set-TextValue 16 4 "123"

# Row 17
Set-TextValue 17 1 "123"
Set-TextValue 17 2 "123"
Set-TextValue 17 4 "123"

# Row 18
Set-TextValue 18 1 "12345"
Set-TextValue 18 2 "1"
Set-TextValue 18 4 "1"

# Rows 19-21: A, B, D populated as text, and column C present but blank
# (copied from the pre-existing blank C8 cell so the empty cell itself,
# not just its value, is carried over/preserved rather than being newly
# assigned and dropped).
Set-TextValue 19 1 "123"
Set-TextValue 19 2 "123"
Set-TextValue 19 4 "123"
$ws.Range("C8").Copy($ws.Range("C19"))

Set-TextValue 20 1 "123"
Set-TextValue 20 2 "123"
Set-TextValue 20 4 "123"
$ws.Range("C8").Copy($ws.Range("C20"))

Set-TextValue 21 1 "123"
Set-TextValue 21 2 "123"
Set-TextValue 21 4 "123"
$ws.Range("C8").Copy($ws.Range("C21"))

# --- Row 8: the previously blank C8 cell is removed entirely ---
$ws.Cells.Item(8, 3).Value = ""
